# Runtime update: refresh "as_of_utc" timestamps for all referee rows on the
# "Главные" and "Линейные" sheets, and bump the underlying stats for the
# referees whose games-played counters increased since the last run.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-10-28 03:08:23"

# Column letter -> index map for readability (A=1 ... AA=27)
$colC = 3
$colD = 4
$colE = 5
$colF = 6
$colG = 7
$colH = 8
$colI = 9
$colJ = 10
$colK = 11
$colV = 22
$colW = 23
$colAA = 27

# Per-sheet, per-row overrides of the stat columns (C,D,E,F,G,H,I,J,K,V,W).
# Only the rows/columns that actually changed are listed; anything not
# listed here keeps its original value and just gets its timestamp bumped.
$updates = @{
    "Главные" = @{
        5  = @{ C = 18; D = 303; E = 161; F = 142; G = 16.83; H = 8.94; I = 7.89; J = 78; K = 71; V = 14; W = 14 }
        7  = @{ C = 13; D = 177; E = 74;  F = 103; G = 13.62; H = 5.69; I = 7.92; J = 37; K = 34; V = 10; W = 4 }
        20 = @{ C = 17; D = 291; E = 112; F = 179; G = 17.12; H = 6.59; I = 10.53; J = 51; K = 62 }
        25 = @{ C = 18; D = 332; E = 172; F = 160; G = 18.44; H = 9.56; I = 8.89; J = 81; K = 75 }
    }
    "Линейные" = @{
        8  = @{ C = 16; D = 244; E = 98;  F = 146; G = 15.25; H = 6.13; I = 9.13; J = 44; K = 58 }
        9  = @{ C = 17; D = 306; E = 142; F = 164; G = 18;    H = 8.35; I = 9.65; J = 66; K = 77 }
        20 = @{ C = 14; D = 229; E = 108; F = 121; G = 16.36; H = 7.71; I = 8.64; J = 54; K = 58; V = 14; W = 16 }
        26 = @{ C = 18; D = 264; E = 129; F = 135; G = 14.67; H = 7.17; I = 7.5;  J = 62; K = 55; V = 14; W = 16 }
    }
}

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowUpdates = $updates[$sheetName]

    for ($r = 2; $r -le 26; $r++) {
        if ($rowUpdates.ContainsKey($r)) {
            $cellChanges = $rowUpdates[$r]
            if ($cellChanges.ContainsKey("C")) { $ws.Cells.Item($r, $colC).Value = $cellChanges["C"] }
            if ($cellChanges.ContainsKey("D")) { $ws.Cells.Item($r, $colD).Value = $cellChanges["D"] }
            if ($cellChanges.ContainsKey("E")) { $ws.Cells.Item($r, $colE).Value = $cellChanges["E"] }
            if ($cellChanges.ContainsKey("F")) { $ws.Cells.Item($r, $colF).Value = $cellChanges["F"] }
            if ($cellChanges.ContainsKey("G")) { $ws.Cells.Item($r, $colG).Value = $cellChanges["G"] }
            if ($cellChanges.ContainsKey("H")) { $ws.Cells.Item($r, $colH).Value = $cellChanges["H"] }
            if ($cellChanges.ContainsKey("I")) { $ws.Cells.Item($r, $colI).Value = $cellChanges["I"] }
            if ($cellChanges.ContainsKey("J")) { $ws.Cells.Item($r, $colJ).Value = $cellChanges["J"] }
            if ($cellChanges.ContainsKey("K")) { $ws.Cells.Item($r, $colK).Value = $cellChanges["K"] }
            if ($cellChanges.ContainsKey("V")) { $ws.Cells.Item($r, $colV).Value = $cellChanges["V"] }
            if ($cellChanges.ContainsKey("W")) { $ws.Cells.Item($r, $colW).Value = $cellChanges["W"] }
        }

        # Every row's "as_of_utc" timestamp is refreshed to the new run time.
        $ws.Cells.Item($r, $colAA).Value = $newTimestamp
    }
}
